$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (65) for phone number 51616170 with 0 points.
# Column A holds the phone number as text (matching the sheet's existing
# text-stored "phone" values), column B (birthday) is left blank/empty
# text, and column C (total_points) is a numeric 0.

# Leading "'" forces text storage (keeps the value looking like a number
# but stored as a string) without leaving a stray quotePrefix style on
# the cell once we reset it back to the default "Normal" style below.
$ws.Range("A65").Value = "'51616170"
$ws.Range("A65").Style = "Normal"

$ws.Range("B65").Value = "'"
$ws.Range("B65").Style = "Normal"

$ws.Range("C65").Value = 0
